$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C, rows 2 through 252 all currently hold 7573; update them to 7293.
$ws.Range("C2:C252").Value = 7293
